# Applies the OOXML changes described by the commit "Solucionado error de tipeo".
#
# Content-level changes (the xmlns/mc:AlternateContent churn visible in the raw
# XML diff is just an artifact of the authoring tool's re-serialization and is
# not something the PowerPoint object model lets us control, so it is not
# reproduced here):
#
#  1. Slide 10: picture "Content Placeholder 4" moved down slightly (Top).
#  2. Slide 15: caption textbox narrowed (Width) and its text corrected/split
#     from "Fig. 4: modelo de datos de la base de datos CMH." into
#     "Fig. 4: Modelo de datos de CMH." (typo fix -> matches commit message).
#  3. Slide 19: caption runs "Fig. " / "6: Tareas de la iteraci" /
#     "ón 2 basada en Gantt." consolidated into a single run.
#  4. Slide 2: agenda bullet runs "Capa de negocios terminal (en proceso" / ")"
#     consolidated into a single run.
#  5. Slide 6: bullet runs for the "Web app" and "Payment service" bullets each
#     consolidated into a single run.
#  6. Slide 7: stray trailing endParaRPr removed from the RFC bullet.

$p = $ppt.ActivePresentation

# --- 1. Slide 10: nudge the picture down (Top: 148.2466pt -> 151.90756pt) ---
$s10 = $p.Slides.Item(10)
$pic10 = $s10.Shapes.Item(3)
$pic10.Top = 151.90756

# --- 2. Slide 15: narrow the caption textbox and fix/split its text ---
$s15 = $p.Slides.Item(15)
$cap15 = $s15.Shapes.Item(4)
$cap15.Width = 287.93465

$tr15 = $cap15.TextFrame.TextRange
$tr15.Text = "placeholder"
$tr15.Text = "Fig. 4: "
[void]$tr15.InsertAfter("Modelo ")
[void]$tr15.InsertAfter("de datos ")
[void]$tr15.InsertAfter("de ")
[void]$tr15.InsertAfter("CMH.")

# --- 3. Slide 19: merge the Gantt caption runs into one ---
$s19 = $p.Slides.Item(19)
$cap19 = $s19.Shapes.Item(3)
$tr19 = $cap19.TextFrame.TextRange
$tr19.Text = "placeholder"
$tr19.Text = "Fig. 6: Tareas de la iteración 2 basada en Gantt."

# --- 4. Slide 2: merge the "Capa de negocios terminal (en proceso)" runs ---
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2)
$para2 = $body2.TextFrame.TextRange.Paragraphs(6)
$para2.Text = "placeholder"
$para2.Text = "Capa de negocios terminal (en proceso)"

# --- 5. Slide 6: merge the "Web app" and "Payment service" bullet runs ---
$s6 = $p.Slides.Item(6)
$body6 = $s6.Shapes.Item(3)
$tr6 = $body6.TextFrame.TextRange

$paraWeb = $tr6.Paragraphs(2)
$paraWeb.Text = "placeholder"
$paraWeb.Text = "Web app: Sitio web que provee las funcionalidades para los pacientes."

$paraPay = $tr6.Paragraphs(3)
$paraPay.Text = "placeholder"
$paraPay.Text = "Payment service: Servicio automatizado de pagos de honorarios."

# --- 6. Slide 7: drop the stray trailing endParaRPr after the RFC bullet ---
$s7 = $p.Slides.Item(7)
$body7 = $s7.Shapes.Item(3)
$tf7 = $body7.TextFrame
$tf7.DeleteText()
$tf7.TextRange.Text = "Creación de documento RFC debido al cambio en la arquitectura."
